# Apply the cryptos list refresh (GitHub Actions data update).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.402.09'
$ws.Range('E2').Value = '  +0.63%  '
$ws.Range('D3').Value = '1.606.95'
$ws.Range('E4').Value = '  -0.01%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range('D5').Value = '212.35'
$ws.Range('E5').Value = '  -0.02%  '
$ws.Range('E6').Value = '  -0.26%  '
$ws.Range('E7').Value = '  +0.02%  '
$ws.Range('E8').Value = '  -0.14%  '
$ws.Range('E9').Value = '  +0.00%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range('D10').Value = '19.28'
$ws.Range('E10').Value = '  +1.38%  '
$ws.Range('E11').Value = '  +0.57%  '
$ws.Range('D12').Value = '1.833.67'
$ws.Range('D13').Value = '1.614.20'
$ws.Range('E13').Value = '  +1.37%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range('D15').Value = '0.507'
$ws.Range('E15').Value = '  -0.32%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range('D16').Value = '63.46'
$ws.Range('E16').Value = '  -0.53%  '
$ws.Range('B17').Value = 'WrappedBTC'
$ws.Range('C17').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range('D17').Value = '26.407.63'
$ws.Range('E17').Value = '  +0.62%  '
$ws.Range('B18').Value = 'BitcoinCash'
$ws.Range('C18').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D18").NumberFormat = "@"
$ws.Range('D18').Value = '232.67'
$ws.Range('E18').Value = '  +7.97%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range('D19').Value = '7.69'
$ws.Range('E19').Value = '  +5.06%  '
$ws.Range('D20').Value = '0.0₃0726'
$ws.Range('E20').Value = '  -0.20%  '
$ws.Range('E21').Value = '  +0.06%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range('D22').Value = '4.27'
$ws.Range('E22').Value = '  -0.41%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range('D23').Value = '8.96'
$ws.Range('E23').Value = '  -0.88%  '
$ws.Range('E24').Value = '  +1.31%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range('D25').Value = '147.22'
$ws.Range('E25').Value = '  +1.93%  '
$ws.Range('E26').Value = '  +0.04%  '
$ws.Range('E27').Value = '  +0.26%  '
$ws.Range('E28').Value = '  +1.18%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range('D29').Value = '15.46'
$ws.Range('E29').Value = '  +2.32%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range('D30').Value = '0.0495'
$ws.Range('E30').Value = '  +1.08%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range('D31').Value = '1.15'
$ws.Range('E31').Value = '  +0.02%  '
$ws.Range('D32').Value = '1.488.55'
$ws.Range('E32').Value = '  +5.37%  '
$ws.Range('E33').Value = '  +1.33%  '
$ws.Range('E34').Value = '  -0.43%  '
$ws.Range('E36').Value = '  +0.98%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range('D37').Value = '0.562'
$ws.Range('E37').Value = '  -3.22%  '
$ws.Range('E38').Value = '  -0.21%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range('D39').Value = '0.822'
$ws.Range('E39').Value = '  -0.11%  '
$ws.Range('E40').Value = '  -0.55%  '
$ws.Range('E41').Value = '  +0.13%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range('D42').Value = '2.19'
$ws.Range('E42').Value = '  +2.49%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range('D43').Value = '0.936'
$ws.Range('E43').Value = '  -3.99%  '
$ws.Range('D44').Value = '1.745.25'
$ws.Range('E44').Value = '  +0.97%  '
$ws.Range('E45').Value = '  -0.59%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range('D46').Value = '60.91'
$ws.Range('E46').Value = '  +0.10%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range('D47').Value = '89.33'
$ws.Range('E47').Value = '  +3.69%  '
$ws.Range('E48').Value = '  +0.46%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range('D49').Value = '0.0501'
$ws.Range('E49').Value = '  +0.08%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range('D50').Value = '0.0962'
$ws.Range('E50').Value = '  +0.85%  '
$ws.Range('B51').Value = 'EnergySwap'
$ws.Range('C51').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D51").NumberFormat = "@"
$ws.Range('D51').Value = '7.44'
$ws.Range('E51').Value = '  +1.43%  '
